# Add data for 2022-08-15: the month-to-date column (column B, "August 2022")
# now covers through August 07 instead of August 06, and new carjacking
# incidents have been tallied in both the current month-to-date column and
# several historical "August" columns (back-filled records) across various
# neighborhood rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and update the column header text to reflect the
# new as-of date.
$ws.Name = "Through 2022-08-07"
$ws.Range("B1").Value = "August 2022 (through August 07)"

# --- Updates to existing values ---
$ws.Range("B2").Value = 5      # Austin / August 2022
$ws.Range("BF2").Value = 3     # Austin / August 2015
$ws.Range("J5").Value = 7      # Garfield Park / August 2021
$ws.Range("R8").Value = 2      # South Shore / August 2020
$ws.Range("J10").Value = 2     # Little Village / August 2021
$ws.Range("R12").Value = 4     # West Town / August 2020
$ws.Range("J13").Value = 2     # Roseland / August 2021
$ws.Range("BF15").Value = 3    # Auburn Gresham / August 2015

# --- New values in previously-empty cells ---
$ws.Range("B4").Value = 1      # North Lawndale / August 2022
$ws.Range("AX6").Value = 1     # Englewood / August 2016
$ws.Range("AP8").Value = 2     # South Shore / August 2017
$ws.Range("Z14").Value = 1     # Woodlawn / August 2019
$ws.Range("Z15").Value = 1     # Auburn Gresham / August 2019
$ws.Range("J22").Value = 1     # Loop / August 2021
$ws.Range("B26").Value = 1     # River North / August 2022
$ws.Range("AP31").Value = 1    # West Loop / August 2017
$ws.Range("AX34").Value = 1    # Washington Park / August 2016
$ws.Range("B35").Value = 1     # West Elsdon / August 2022
$ws.Range("AH49").Value = 1    # Galewood / August 2018
$ws.Range("B50").Value = 1     # Gage Park / August 2022
$ws.Range("J56").Value = 1     # Wrigleyville / August 2021
$ws.Range("Z61").Value = 1     # Avondale / August 2019
$ws.Range("B78").Value = 1     # Lincoln Square / August 2022
$ws.Range("AH95").Value = 1    # Streeterville / August 2018
$ws.Range("B96").Value = 1     # Uptown / August 2022
